$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (mirrors the source diff).
$updates = [ordered]@{
    "D2"  = "59.853.92"
    "E2"  = "  +1.05%  "
    "D3"  = "2.302.58"
    "E3"  = "  -0.62%  "
    "E4"  = "  -0.01%  "
    "D5"  = "542.09"
    "E5"  = "  +0.08%  "
    "D6"  = "129.38"
    "E6"  = "  -2.52%  "
    "D8"  = "0.572"
    "E8"  = "  -2.58%  "
    "D9"  = "2.300.03"
    "E9"  = "  -0.63%  "
    "E10" = "  -0.65%  "
    "D11" = "5.52"
    "E11" = "  +0.87%  "
    "E12" = "  -0.20%  "
    "E13" = "  -0.90%  "
    "D14" = "23.30"
    "E14" = "  -2.90%  "
    "D15" = "59.820.79"
    "E15" = "  +1.22%  "
    "D16" = "2.710.80"
    "E16" = "  -0.66%  "
    "E17" = "  -1.38%  "
    "D18" = "2.306.83"
    "E18" = "  -0.37%  "
    "D19" = "10.46"
    "E19" = "  -1.78%  "
    "E20" = "  -2.71%  "
    "D21" = "311.28"
    "E21" = "  -0.69%  "
    "E22" = "  -0.78%  "
    "E23" = "  -0.11%  "
    "E24" = "  +0.02%  "
    "D25" = "63.63"
    "E25" = "  +1.54%  "
    "E26" = "  -1.87%  "
    "E27" = "  +0.08%  "
    "D28" = "7.73"
    "E28" = "  -2.98%  "
    "D30" = "1.18"
    "E30" = "  +0.15%  "
    "D31" = "170.28"
    "E31" = "  -0.05%  "
    "E32" = "  -0.83%  "
    "D33" = "0.0₃0724"
    "E33" = "  -2.51%  "
    "D34" = "5.81"
    "E34" = "  -1.48%  "
    "E35" = "  +1.78%  "
    "E36" = "  -2.11%  "
    "E37" = "  +0.02%  "
    "D38" = "17.64"
    "E38" = "  -1.29%  "
    "E39" = "  +0.06%  "
    "E40" = "  -2.72%  "
    "D41" = "317.85"
    "E41" = "  +3.94%  "
    "D42" = "37.98"
    "E42" = "  -0.99%  "
    "E43" = "  -1.40%  "
    "D44" = "135.73"
    "E44" = "  -3.86%  "
    "E45" = "  -1.34%  "
    "E46" = "  -2.59%  "
    "E47" = "  +0.73%  "
    "E48" = "  +1.09%  "
    "D49" = "0.0488"
    "E49" = "  -1.66%  "
    "E50" = "  +19.29%  "
    "E51" = "  -0.57%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Keep these as plain text, matching the workbook's existing inline-string
    # cells (prices with multiple dots, padded percentage strings, etc.).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
